$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.407.27'
$ws.Range('E2').Value = '  +2.84%  '
$ws.Range('D3').Value = '2.306.20'
$ws.Range('E3').Value = '  +1.74%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '310.78'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.38'
$ws.Range('E6').Value = '  +6.15%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.532'
$ws.Range('E7').Value = '  +1.21%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.527'
$ws.Range('E9').Value = '  +7.76%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.61'
$ws.Range('E10').Value = '  +4.18%  '
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '52.14'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('E13').Value = '  -1.15%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.03'
$ws.Range('E14').Value = '  +3.24%  '
$ws.Range('D15').Value = '2.663.58'
$ws.Range('E15').Value = '  +1.75%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.07'
$ws.Range('E16').Value = '  +3.28%  '
$ws.Range('D17').Value = '2.315.89'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.811'
$ws.Range('E18').Value = '  +2.55%  '
$ws.Range('D19').Value = '43.287.71'
$ws.Range('E19').Value = '  +3.00%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.20'
$ws.Range('E20').Value = '  -0.45%  '
$ws.Range('E21').Value = '  +2.21%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.17'
$ws.Range('E22').Value = '  +3.52%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '68.15'
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '243.04'
$ws.Range('E24').Value = '  +2.75%  '
$ws.Range('E25').Value = '  +2.88%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.61'
$ws.Range('E26').Value = '  +1.09%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '24.83'
$ws.Range('E28').Value = '  +5.43%  '
$ws.Range('E29').Value = '  +8.12%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '37.02'
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '9.65'
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '167.58'
$ws.Range('E32').Value = '  +2.12%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.28'
$ws.Range('E33').Value = '  +0.73%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '18.19'
$ws.Range('E35').Value = '  +4.93%  '
$ws.Range('E36').Value = '  +6.13%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0745'
$ws.Range('E38').Value = '  -2.18%  '
$ws.Range('E39').Value = '  +3.20%  '
$ws.Range('E40').Value = '  +2.05%  '
$ws.Range('E41').Value = '  +7.76%  '
$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.79'
$ws.Range('E42').Value = '  +22.24%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.116'
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0294'
$ws.Range('E44').Value = '  +4.64%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.979.18'
$ws.Range('E45').Value = '  +0.97%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '19.06'
$ws.Range('E46').Value = '  +1.25%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.00'
$ws.Range('E47').Value = '  +2.82%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.99'
$ws.Range('E48').Value = '  +1.93%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '55.98'
$ws.Range('E49').Value = '  +4.49%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.93'
$ws.Range('E50').Value = '  +1.22%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.59'
$ws.Range('E51').Value = '  +8.31%  '
